$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "fh5ZTeZxj1K8JDT4"
$ws.Range("A3").Value = "a32y5B8fUezXbwsZ"
$ws.Range("A4").Value = "eX1dEoJRscW6q3JY"
